$d = $word.ActiveDocument

# 1. Update "Curso (semestre ideal)" line to add "EF (7), "
$d.Content.Find.Execute("Curso (semestre ideal): EP (3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EF (7), EP (3)", 2)

# 2. Fix accent in "Fabrício" -> "Fabricio"
$d.Content.Find.Execute("5840917 - Fabrício Maciel Gomes", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5840917 - Fabricio Maciel Gomes", 2)

# 3. Remove the "Requisitos" heading paragraph and the "LOQ4203..." paragraph entirely.
$found = $true
while ($found) {
    $found = $false
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t -like "Requisitos*" -or $t -like "LOQ4203*") {
            $p.Range.Delete()
            $found = $true
            break
        }
    }
}
